$d = $word.ActiveDocument

# 1. Remove the paragraph that holds the two inline drawings (the picture-only
#    paragraph right after old question (3)). Find it by InlineShapes.Count > 0
#    rather than a hard-coded index, so the scan is robust either way.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.InlineShapes.Count -gt 0) {
        $para.Range.Delete()
    }
}

# 2. Title
$d.Paragraphs.Item(1).Range.Text = '盈虧問題 - 基本型'

# 3. Replace questions (1)-(7) in place
$d.Paragraphs.Item(2).Range.Text = '(1) 老猴子給小猴子分梨。每隻小猴子分6個梨，多出12個梨；每隻小猴子分7個梨，少11個梨。問：有【  】隻小猴子、【   】個梨'
$d.Paragraphs.Item(3).Range.Text = '(2) 重陽節那天，六零一班的少先隊員帶了一些蘋果去敬老院慰問老人。如果每人分11個，則剩下39個；如果每人分14個，則剩下12個。問：有老人【 】個，蘋果有【   】個'
$d.Paragraphs.Item(4).Range.Text = '(3) 學校買來一批圖書。若每人發9本，則少25本；若每人發6本，則少7本。問：有【 】個學生，買了【  】本圖書'
$d.Paragraphs.Item(5).Range.Text = '(4) 若「甲數的13倍減6」與「甲數的7倍加18」相等，則甲數的8倍減12＝【  】'
$d.Paragraphs.Item(6).Range.Text = '(5) 將一批本子發給學生，每人發10本，差90本；若每人發8本，則仍差8本。問：學生【  】人和本子【   】本'
$d.Paragraphs.Item(7).Range.Text = '(6) 一批小朋友去買東西，若每人出10元則多8元；若每人出7元則少4元。問：有【 】個小朋友，東西的價格是【  】元'
$d.Paragraphs.Item(8).Range.Text = '(7) 同學們植樹，如果每人種2棵，還有18棵沒種；如果每個人種5棵，還有3顆沒有種。問：有【 】個同學，有【  】棵樹'

# 4. Append questions (8)-(40) as new paragraphs after the last one,
#    cloning the "question" style from the paragraph they split off from.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(8) 王老師給美術社的同學分發圖畫紙。如果每人發5張，則少32張；如果每人發3張，則少2張。問：美術社有【  】名同學，王老師一共有【  】張圖畫紙'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(9) 301班少先隊員參加搬磚勞動。如果每人搬4塊磚，還剩7塊；如果每人搬5塊，則少2塊磚。問：這個班少先隊有【 】個人、要搬的磚共有【  】塊'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(10) 一個植樹小組植樹。如果每人栽5棵，還剩14棵；如果每人栽7棵，就缺4棵。問：這個植樹小組有【 】人、一共有【  】棵樹'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(11) 小明的爺爺買回一筐梨，分給全家人。如果小明和小妹每人分4個梨，其餘每人分2個梨，還多出4個梨。如果小明1人分6個梨，其餘每人分4個梨，又差12個梨。小明家有【 】人、這筐梨子有【  】個'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(12) 同學們暑假前到圖書館借書，如果每人借4本，則最後少2本；如果前2人每人先借8本，餘下的人每人借3本，這些圖書恰好借完。問：書的總數是【  】本'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(13) 工人鋪一條路基，若每天鋪260公尺，鋪完全路長就得延長8天；若每天鋪300公尺，鋪完全路長仍要延長4天。問：這條路長【    】公尺。'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(14) 老師分配若干人擦玻璃。如果其中2人各擦4塊，其餘每人擦5塊，則餘22塊；如果每人擦7塊，正好擦完。問：擦玻璃的有【  】人，玻璃【  】塊'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(15) 築路隊計畫每天築路720公尺，實際每天比原計劃多築80公尺，這樣在完成規定任務的前三天，就只剩下1160公尺未築。問：這條路長【     】公尺'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(16) 三年級給優秀學生發獎品書，如果每個學生發5冊還剩32冊；如果其中10個學生每人發4冊，其餘每人發8冊，就恰好發完。問：優秀學生有【  】人、獎品書有【   】冊'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(17) 全班同學分組勞動，每組8人。勞動中覺得每組人數太少，因而重新編組，每組改為12人，這樣減少了2組，問：參加勞動的學生有【  】人'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(18) 學校買來一批電風扇分給各班。若有兩個班每班分到4台，其餘每班只能分2台；如果有一個班分6台，其餘每班分4台，還差12台。問：共買來【  】台電風扇、有【 】個班'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(19) 王師傅加工一批零件，每天加工20個，可以提前1天完成。工作4天後，由於改進了技術，每天可多加工5個，結果提前3天完成。問：這批零件有【   】個'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(20) 王老師看一本小說，如果每天看25頁，看完全書比規定時間多一天，如果每天看30頁最後一天少看15頁，如果每天看29頁，最後一天要讀【  】頁才能按規定讀完'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(21) 一個旅遊團去旅館住宿，6人一間，多2個房間；若4人一間又少2個房間。問：旅遊團共有【  】人。'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(22) 學生參加校外教學，如果每輛車坐30個學生，則多12個座位；如果每輛車坐36個學生，則多兩輛車。問：有【  】輛車，【   】個學生'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(23) 學校春遊租了幾條船讓學生們劃，如果每條船坐3人，則多出一條船；如果每條船坐5人，則空出19人的位置，問：有學生【  】人參加划船。'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(24) 若干個同學去划船，他們租了一些船。若每船坐4人，則多5人；若每船坐5人，則船上有4個空位。問，這個小組一共有【  】同學、租了【 】條船。'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(25) 五年級同學去划船，如果增加一隻船，正好每只船上坐7人；如果減少一隻船，正好每只船上坐8人。問：這個年級共有【   】個同學。'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(26) 同學去划船，如果每只船坐4人，則少3只船；如果每只船坐6人，則少2人，問：同學們共【  】人、租了【 】隻船。'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(27) 同學們去公園划船，每條船坐4人，就會少3條船；每條船坐6人，還有2人坐不下。問：一共有【  】個同學、小船有【 】條'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(28) 學校為新生分配宿舍．每個房間住3人，則多出22人；每個房間多住5人，則空1個房間。問：宿舍有【 】間、新生有【  】人。'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(29) 四年三班的一部分同學去野餐。如果每張餐布周圍坐4名同學，就有6名同學沒有座位；如果每張餐布周圍多坐1名同學，就會多出4個位置。問：有【  】名同學參加野餐、帶了【  】張餐布'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(30) 陽光小學學生乘遊覽車到陽明山郊遊。如果每車坐65人，則有5人不能乘上車；如果每車多坐5人，恰多餘了一輛車。問：共有【  】輛汽車，有【   】名學生'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(31) 用繩子測量井深。如果把繩子三折垂到水面，餘7公尺；如果把繩子5折垂到水面，餘1公尺。問：繩長【  】公尺、井深【 】公尺'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(32) 用一根繩子測井口到水面的高度，把繩子對折後垂直到水面，繩子超過井口9公尺；   把繩子三折後垂直到水面，繩子超過井口3公尺。問：繩子長【  】公尺、井口到水面的距離【 】公尺'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(33) 用繩子測游泳池的深度，把繩子對折測量時，池外多出4公尺，把繩子三折測量時，繩子離池口還差2公尺。問：井深【  】公尺，繩子長【  】公尺'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(34) 小華用自己僅存的漆包線在磁棒上繞線圈，當他繞了80圈時，測得餘線長15.28   公分，於是想改繞90圈，卻發現缺少22.4公分的漆包線。問：王華的漆包線長【      】公分，所用的磁棒54的半徑【   】公分'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(35) 用一根繩子測井深，把繩子折3折來測，則井外餘1公尺，把繩子折5折來測，則繩子離井口還差1公尺。問：井深【 】公尺、繩子長【  】公尺'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(36) 用繩子測井深，把繩子二折來量，井外餘5公尺；把繩子三折來量，還差1公尺。問：井深【  】公尺，繩子長【  】公尺'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(37) 用繩子測井深，將繩子兩折來測，則井外餘80公分；將繩子3折來測，則繩子還差70公分。問：井深【   】公分，繩子長【   】公分'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(38) 用一根繩子測量壩長，用整根繩子量30下，堤壩還剩4公尺，用繩子的一半量56下，堤壩還剩16公尺，求：繩子長度【 】公尺、堤壩長【   】公尺'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(39) 用一根繩繞樹5周還剩公尺，若用繩的三分之一繞樹一周還餘公尺。問：繩長【 】公尺，樹的周長【】公尺'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = '(40) 用繩子測游泳池水深，繩子兩折時，多餘60公分；繩子三折時，還差40公分。求：繩長【   】公分、游泳池水深【   】公分'

